# This workbook's "Artfynd" sheet contains 15 data rows (rows 2-16, with row 1
# being the header). The edit is a permutation of those 15 rows: each
# destination row ends up holding the full record that originally lived at a
# different row. Columns outside A:AY / data columns are untouched (they are
# identical across all rows anyway).
#
# Strategy: stage a full copy of each source row (2-16) into a scratch area
# far below the used range, then copy each staged row back into its final
# destination row. Using Range.Copy (rather than re-typing cell Values)
# preserves the original cell data types (e.g. numeric-looking values stored
# as text) without introducing spurious styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 16
$lastCol = "AY"

# Row in current sheet -> row in scratch staging area (offset well below data)
$stagingOffset = 1000

# Step 1: copy each source row into its staging slot
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRange = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $stagingRow = $r + $stagingOffset
    $dstRange = $ws.Range("A" + $stagingRow + ":" + $lastCol + $stagingRow)
    $srcRange.Copy($dstRange)
}

# Step 2: destination row -> original source row (before the edit) mapping
$mapping = @{
    2  = 16
    3  = 11
    4  = 9
    5  = 8
    6  = 10
    7  = 13
    8  = 7
    9  = 3
    10 = 5
    11 = 2
    12 = 6
    13 = 12
    14 = 15
    15 = 4
    16 = 14
}

# Step 3: copy each staged row back into its new destination row
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $stagingRow = $srcRow + $stagingOffset
    $srcRange = $ws.Range("A" + $stagingRow + ":" + $lastCol + $stagingRow)
    $dstRange = $ws.Range("A" + $destRow + ":" + $lastCol + $destRow)
    $srcRange.Copy($dstRange)
}

# Step 4: clear out the scratch staging area
$clearRange = $ws.Range("A" + (1000 + $firstRow) + ":" + $lastCol + (1000 + $lastRow))
$clearRange.Clear()

# Step 5: every data row only ever populates this specific set of columns;
# copying the full A:AY block (above) also materializes empty cells for the
# columns that are never used, so strip those back out to match the
# original sparse layout.
$neverUsedCols = @('L','M','O','X','AC','AH','AI','AJ','AK','AL','AM','AN','AO','AP','AQ','AR','AS','AU','AV')
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $neverUsedCols) {
        $ws.Range($c + $r).Clear()
    }
}
